$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D "TODO" notes (independent running list in column D) ---
# These two are written first so the shared-string table picks up the
# same index ordering the original authoring tool produced.
$ws.Range("D61").Value = "fix: tfidf on server"
$ws.Range("D62").Value = "fix: connect to db on vscode"

# --- Row 57: new Task-log entry (columns A/B) ---
# Clone formatting from the row above (A56/B56) so the date- and
# task-column styling (borders etc.) matches the rest of the log,
# then overwrite with the new values.
$ws.Range("A56").Copy()
$ws.Range("A57").PasteSpecial(-4122)
$ws.Range("B56").Copy()
$ws.Range("B57").PasteSpecial(-4122)
$ws.Range("A57").Value = 45204
$ws.Range("B57").Value = "2971 docs inserted on server db, BA: topic modeling init"
$ws.Rows.Item(57).RowHeight = 35

# --- Row 58: new Task-log entry (columns A/B) ---
$ws.Range("A56").Copy()
$ws.Range("A58").PasteSpecial(-4122)
$ws.Range("B56").Copy()
$ws.Range("B58").PasteSpecial(-4122)
$ws.Range("A58").Value = 45205
$ws.Range("B58").Value = "PCA & cluster in DB"
$ws.Rows.Item(58).RowHeight = 18

$ws.Application.CutCopyMode = 0

# --- More column D "TODO" notes, continuing after the row-63 gap ---
$ws.Range("D64").Value = "BA: PCA dim update"
$ws.Range("D65").Value = "BERTopic"
$ws.Range("D66").Value = "LDA"

# --- Row 60: stray single-space note in column B ---
$ws.Range("B60").Value = " "

# Leave the sheet scrolled/selected near the new last entry, like the
# author's Excel window would have been positioned after typing it.
$ws.Range("D64").Select()
